$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 5291
$ws.Range("F6").Value = 5291
$ws.Range("F7").Value = 160
$ws.Range("F10").Value = 6
$ws.Range("F11").Value = 1190
$ws.Range("F12").Value = 749
$ws.Range("F13").Value = 5229
$ws.Range("F14").Value = 32
$ws.Range("F15").Value = 75
$ws.Range("F16").Value = 93
$ws.Range("F17").Value = 2397
$ws.Range("F18").Value = 2397
$ws.Range("F19").Value = 259
$ws.Range("F20").Value = 103
$ws.Range("F22").Value = 3935
$ws.Range("F26").Value = 3864
$ws.Range("F28").Value = 179
$ws.Range("F29").Value = 251
$ws.Range("F30").Value = 211
$ws.Range("F37").Value = 6850
$ws.Range("F38").Value = 1115
$ws.Range("F39").Value = 535
$ws.Range("F42").Value = 1396
$ws.Range("F44").Value = 706
$ws.Range("F46").Value = 2334
$ws.Range("F49").Value = 10
$ws.Range("F50").Value = 788
$ws.Range("F51").Value = 932

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 90
$ws.Range("F14").Value = 7
$ws.Range("F17").Value = 143
$ws.Range("F25").Value = 822

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 217

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 217
$ws.Range("F7").Value = 5291
$ws.Range("F8").Value = 5291
$ws.Range("F9").Value = 160
$ws.Range("F11").Value = 90
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 1190
$ws.Range("F15").Value = 749
$ws.Range("F16").Value = 32
$ws.Range("F17").Value = 75
$ws.Range("F18").Value = 93
$ws.Range("F19").Value = 2397
$ws.Range("F20").Value = 2398
$ws.Range("F21").Value = 259
$ws.Range("F22").Value = 103
$ws.Range("F24").Value = 3935
$ws.Range("F25").Value = 3864
$ws.Range("F27").Value = 179
$ws.Range("F28").Value = 251
$ws.Range("F29").Value = 211
$ws.Range("F35").Value = 143
$ws.Range("F36").Value = 6850
$ws.Range("F37").Value = 1115
$ws.Range("F38").Value = 535
$ws.Range("F42").Value = 1396
$ws.Range("F44").Value = 707
$ws.Range("F46").Value = 2334
$ws.Range("F49").Value = 788
$ws.Range("F50").Value = 932

